$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# Keep the first three headers, but replace the fourth and add a fifth.
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "penjualan_kode"
$ws.Range("C1").Value = "pembeli"
$ws.Range("D1").Value = "barang_id"
$ws.Range("E1").Value = "jumlah"

# Header formatting: A1:C1 bold + centered, D1:E1 bold only.
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("D1:E1").Font.Bold = $true

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "TXR0015"
$ws.Range("C2").Value = "Customer 15"
$ws.Range("D2").ClearFormats()
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 1

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = 3
$ws.Range("C3").Value = "Customer 16"
$ws.Range("B3").Value = "TXR0016"
$ws.Range("D3").ClearFormats()
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 1

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4:C4").ClearContents()
$ws.Range("D4").ClearFormats()
$ws.Range("D4").Value = 28
$ws.Range("E4").Value = 2

# --- Column widths ---------------------------------------------------
$ws.Range("B:B").ColumnWidth = 20.26953125
$ws.Range("C:C").ColumnWidth = 14.54296875
$ws.Range("D:D").ColumnWidth = 13.7265625

# --- Selection -------------------------------------------------------
$ws.Range("I14").Select()
